$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/27/2023  Through  4/2/2023"

# --- Precinct crime-stat table updates (rows 15-30) ---
$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = -75
$ws.Range("L15").Value = -60
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -18.75
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 47
$ws.Range("K16").Value = -6.382978723404
$ws.Range("L16").Value = 37.5
$ws.Range("M16").Value = -8.333333333333
$ws.Range("N16").Value = -81.893004115226
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 53.333333333333
$ws.Range("I17").Value = 63
$ws.Range("J17").Value = 55
$ws.Range("K17").Value = 14.545454545454
$ws.Range("L17").Value = 31.25
$ws.Range("M17").Value = 142.307692307692
$ws.Range("N17").Value = 8.620689655172
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -85.714285714285
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -55.555555555555
$ws.Range("I18").Value = 51
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = -26.086956521739
$ws.Range("L18").Value = -1.923076923076
$ws.Range("M18").Value = -42.045454545454
$ws.Range("N18").Value = -87.281795511221
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -5.263157894736
$ws.Range("F19").Value = 90
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = 32.352941176470
$ws.Range("I19").Value = 254
$ws.Range("J19").Value = 225
$ws.Range("K19").Value = 12.888888888888
$ws.Range("L19").Value = 68.211920529801
$ws.Range("M19").Value = -28.248587570621
$ws.Range("N19").Value = -61.102603369065
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 3
$ws.Range("I20").Value = 15
$ws.Range("J20").Value = 12
$ws.Range("K20").Value = 25
$ws.Range("L20").Value = 400
$ws.Range("M20").Value = 114.285714285714
$ws.Range("N20").Value = -95.268138801261
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -6.666666666666
$ws.Range("F21").Value = 145
$ws.Range("H21").Value = 9.022556390977
$ws.Range("I21").Value = 429
$ws.Range("J21").Value = 416
$ws.Range("K21").Value = 3.125
$ws.Range("L21").Value = 47.422680412371
$ws.Range("M21").Value = -18.285714285714
$ws.Range("N21").Value = -74.433849821215
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 23
$ws.Range("J22").Value = 29
$ws.Range("K22").Value = -20.689655172413
$ws.Range("L22").Value = 130
$ws.Range("M22").Value = 9.523809523809
$ws.Range("M23").Value = -40
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = -6.976744186046
$ws.Range("F24").Value = 146
$ws.Range("G24").Value = 178
$ws.Range("H24").Value = -17.977528089887
$ws.Range("I24").Value = 484
$ws.Range("J24").Value = 570
$ws.Range("K24").Value = -15.087719298245
$ws.Range("L24").Value = 9.255079006772
$ws.Range("M24").Value = 11.520737327188
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 12.195121951219
$ws.Range("I25").Value = 117
$ws.Range("J25").Value = 119
$ws.Range("K25").Value = -1.680672268907
$ws.Range("L25").Value = 48.101265822784
$ws.Range("M25").Value = 21.875
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 8
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = -27.272727272727
$ws.Range("L26").Value = 33.333333333333
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 15
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 30
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = 30.434782608695
$ws.Range("L27").Value = 42.857142857142
$ws.Range("G30").Value = 3
